$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates -------------------------------------------------
$meta = $wb.Worksheets.Item("Metadata")

# URL: ibm.com -> linuxforhealth.org
$meta.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/subscriberId-encrypted"

# Version: 7.0.0 -> 8.0.0
$meta.Range("B3").Value = "8.0.0"

# Date: 2022-09-08T16:11:15+00:00 -> 2022-11-10T16:00:46+00:00
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"

# Publisher: Alvearie Team -> LinuxForHealth Team
$meta.Range("B9").Value = "LinuxForHealth Team"

# --- Elements sheet updates --------------------------------------------------
$elements = $wb.Worksheets.Item("Elements")

# The "Extension.url" row's Fixed Value mirrors the same URL shown above.
$elements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/subscriberId-encrypted"

# The root "Extension" row's Constraint(s) cell no longer carries the
# ele-1/ext-1 constraint text (that constraint now only applies to the
# child "Extension.extension" row, which already carries it separately).
$elements.Range("AI2").Value = ""
